# Create EPG Test cases
# Insert a new "EpgScreen" worksheet between "hub" and "screenTitles",
# populate it with lookup data, and add a row to "screenTitles".

$wb = $excel.ActiveWorkbook
$hub = $wb.Worksheets.Item("hub")
$screenTitles = $wb.Worksheets.Item("screenTitles")

# --- 1. Create the new sheet right after "hub" --------------------------
$epg = $wb.Worksheets.Add($null, $hub, 1)
$epg.Name = "EpgScreen"

# --- 2. Populate cell values in the same order the strings were first
#        introduced into the shared-string table ------------------------
$epg.Range("O2").Value = "rgba(242, 249, 250, 1)"
$epg.Range("B2").Value = "programTitle"
$epg.Range("M4").Value = "Proximus, ProximusBold"
$epg.Range("L3").Value = "24px"
$epg.Range("O3").Value = "rgba(128, 128, 128, 1)"
$epg.Range("O4").Value = "rgba(255, 255, 0, 1)"

$screenTitles.Range("A5").Value = "EpgChannelScreen"
$screenTitles.Range("B5").Value = "donderdag, 25 mei"

$epg.Range("L1").Value = "font_size"
$epg.Range("M1").Value = "font_family"
$epg.Range("P1").Value = "No_of_Channel"

$epg.Range("A3").Value = "Seniour_groen_Grijs"
$epg.Range("A4").Value = "Seniour_groen_Geel"
$epg.Range("A2").Value = "Standard"
$epg.Range("L2").Value = "20px"
$epg.Range("M8").Value = " "

# --- 3. Remaining header cells (re-use existing shared strings) ---------
$epg.Range("A1").Value = "objectID"
$epg.Range("B1").Value = "name_nl"
$epg.Range("C1").Value = "focused_icon_textline"
$epg.Range("D1").Value = "non_focused_icon"
$epg.Range("E1").Value = "focused_icon_showcase"
$epg.Range("F1").Value = "hubDefaultX"
$epg.Range("G1").Value = "hubDefaultY"
$epg.Range("H1").Value = "width"
$epg.Range("I1").Value = "height"
$epg.Range("J1").Value = "hubFocusedX"
$epg.Range("K1").Value = "hubFocusedY"
$epg.Range("N1").Value = "font-family-showcase"
$epg.Range("O1").Value = "color"

# --- 4. Remaining body cells (re-use existing shared strings) -----------
$epg.Range("M2").Value = "Proximus, ProximusRegular"
$epg.Range("M3").Value = "Proximus, ProximusRegular"

# --- 5. Numeric cells -----------------------------------------------------
$epg.Range("P2").Value = 6
$epg.Range("P3").Value = 4
$epg.Range("P4").Value = 4

# --- 6. Formatting: header row uses the same style as the other sheets' --
$hub.Range("A1").Copy() | Out-Null
$epg.Range("A1:P1").PasteSpecial(-4122) | Out-Null

# --- 7. Formatting: small Consolas note style ("s=3") --------------------
$note = $epg.Range("L2")
$note.Font.Name = "Consolas"
$note.Font.Size = 9
$note.Font.Color = 2236962

$note.Copy() | Out-Null
$epg.Range("M2").PasteSpecial(-4122) | Out-Null
$epg.Range("M3").PasteSpecial(-4122) | Out-Null
$epg.Range("O3").PasteSpecial(-4122) | Out-Null
$epg.Range("M4").PasteSpecial(-4122) | Out-Null
$epg.Range("M8").PasteSpecial(-4122) | Out-Null
$screenTitles.Range("B5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- 8. Column widths / visibility on the new sheet ----------------------
$epg.Columns.Item(1).ColumnWidth = 26.73
$epg.Columns.Item(2).ColumnWidth = 39.17
$epg.Columns.Item(3).ColumnWidth = 15.58
$epg.Columns.Item(12).ColumnWidth = 8.02
$epg.Columns.Item(13).ColumnWidth = 29.88
$epg.Columns.Item(14).ColumnWidth = 0.3
$epg.Columns.Item(15).ColumnWidth = 22.58
$epg.Columns.Item(16).ColumnWidth = 13.17

$epg.Columns.Item(4).ColumnWidth = 16.58
$epg.Columns.Item(5).ColumnWidth = 21.88
$epg.Columns.Item(6).ColumnWidth = 11.3
$epg.Columns.Item(7).ColumnWidth = 11.3
$epg.Columns.Item(8).ColumnWidth = 8.3
$epg.Columns.Item(9).ColumnWidth = 8.3
$epg.Columns.Item(10).ColumnWidth = 12.02
$epg.Columns.Item(11).ColumnWidth = 12.02

$epg.Range("D1:K1").EntireColumn.Hidden = $true

# --- 9. Column width on screenTitles (content widened col A) -------------
$screenTitles.Columns.Item(1).ColumnWidth = 16.9

# --- 10. Selections / active sheet ---------------------------------------
$screenTitles.Range("A2").Select() | Out-Null
$epg.Activate() | Out-Null
$epg.Range("M3").Select() | Out-Null
